$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range('A5').Value = 'nathumang'
$ws.Range('B5').Value = 'nathumang@gmail.com'
$ws.Range('C5').Value = 'a200ecf982d685b7bffc3d9d7eae082eadfdfed1206ae4c452e0e9c3dd4cfdae'
$ws.Range('D5').Value = '[''Germany'', 2, 0]'
$ws.Range('E5').Value = '[''Hungary'', 2, 1]'
$ws.Range('F5').Value = '[''Draw'', 1, 1]'
$ws.Range('G5').Value = '[''Italy'', 1, 0]'
$ws.Range('H5').Value = '[''Netherlands'', 0, 2]'
$ws.Range('I5').Value = '[''Denmark'', 1, 2]'
$ws.Range('J5').Value = '[''England'', 0, 2]'
$ws.Range('K5').Value = '[''Ukraine'', 1, 2]'
$ws.Range('L5').Value = '[''Belgium'', 3, 0]'
$ws.Range('M5').Value = '[''France'', 1, 3]'
$ws.Range('N5').Value = '[''Turkey'', 2, 1]'
$ws.Range('O5').Value = '[''Portugal'', 3, 0]'

# Row 6
$ws.Range('A6').Value = 'kk_queen'
$ws.Range('B6').Value = 'karimkarishma94@gmail.com'
$ws.Range('C6').Value = '7b4a5a80af236549bdd97321ed83593cefcca97e7e1262aa5a96bdf604d82ef5'

# Row 7
$ws.Range('A7').Value = 'Nami'
$ws.Range('C7').Value = '01a9b0c26e5eb24bde6e64834e93a2b9693eeee3ef882edcd04d171b89a26516'
$ws.Range('D7').Value = '[''Germany'', 2, 1]'

# Row 8
$ws.Range('C8').Value = 'e3b0c44298fc1c149afbf4c8996fb92427ae41e4649b934ca495991b7852b855'

# Row 9
$ws.Range('A9').Value = 'prat.d '
$ws.Range('B9').Value = 'pratyush.devliyal@gmail.com'
$ws.Range('C9').Value = 'd3332bd40f9467474b0326549dc6782ee2e0491e209540b6d79fe3f2f73e9194'
$ws.Range('D9').Value = '[''Germany'', 3, 0]'
$ws.Range('E9').Value = '[''Hungary'', 1, 0]'
$ws.Range('F9').Value = '[''Draw'', 1, 1]'
$ws.Range('G9').Value = '[''Italy'', 1, 0]'
$ws.Range('H9').Value = '[''Draw'', 0, 0]'
$ws.Range('I9').Value = '[''Denmark'', 0, 2]'
$ws.Range('J9').Value = '[''England'', 0, 3]'
$ws.Range('K9').Value = '[''Romania'', 1, 0]'
$ws.Range('L9').Value = '[''Belgium'', 2, 1]'
$ws.Range('M9').Value = '[''France'', 0, 2]'
$ws.Range('N9').Value = '[''Draw'', 0, 0]'
$ws.Range('O9').Value = '[''Portugal'', 1, 0]'
$ws.Range('P9').Value = '[''Croatia'', 2, 1]'
$ws.Range('Q9').Value = '[''Germany'', 2, 0]'
$ws.Range('R9').Value = '[''Scotland'', 2, 1]'
$ws.Range('S9').Value = '[''Draw'', 1, 1]'
$ws.Range('T9').Value = '[''England'', 1, 2]'
$ws.Range('U9').Value = '[''Spain'', 2, 1]'
$ws.Range('V9').Value = '[''Slovakia'', 2, 0]'
$ws.Range('W9').Value = '[''Draw'', 0, 0]'
$ws.Range('X9').Value = '[''France'', 0, 1]'
$ws.Range('Y9').Value = '[''Czech Republic'', 0, 1]'
$ws.Range('Z9').Value = '[''Portugal'', 0, 2]'
$ws.Range('AA9').Value = '[''Draw'', 0, 0]'
$ws.Range('AB9').Value = '[''Germany'', 0, 4]'
$ws.Range('AC9').Value = '[''Hungary'', 2, 3]'
$ws.Range('AD9').Value = '[''Spain'', 1, 3]'
$ws.Range('AE9').Value = '[''Draw'', 1, 1]'
$ws.Range('AF9').Value = '[''Netherlands'', 1, 0]'
$ws.Range('AG9').Value = '[''France'', 2, 1]'
$ws.Range('AH9').Value = '[''England'', 4, 1]'
$ws.Range('AI9').Value = '[''Denmark'', 2, 0]'
$ws.Range('AJ9').Value = '[''Romania'', 0, 1]'
$ws.Range('AL9').Value = '[''Portugal'', 0, 3]'
$ws.Range('AM9').Value = '[''Draw'', 0, 0]'

# Row 10
$ws.Range('A10').Value = 'sumit'
$ws.Range('B10').Value = 'sumitsinhaiitkgp@gmail.com'
$ws.Range('C10').Value = '27891eae1d9194dafa403e168161ae859340e9f6c0102ee3ca0f867669cce9b0'
